# Updates cryptos list price/volume figures (and the ARBITRUM/MXToken row swap)
# per the source commit. Price-column (D) values that look like plain decimal
# numbers (e.g. "220.55") would otherwise be auto-coerced by Excel's Value
# setter into a floating-point number (losing the exact text / adding FP noise
# and a numFmt style). We force text by prefixing a quote (classic "text entry"
# trick), then reset Style back to Normal so no stray quotePrefix/numFmt survives
# on the cell - matching the original workbook where these are plain inline
# strings with no explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'" + '26.335.07'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +0.52%  '

# Row 3
$ws.Cells.Item(3, 4).Value = "'" + '1.668.41'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.76%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'" + '220.55'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.09%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'" + '0.5313'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.03%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'" + '0.2652'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +0.98%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'" + '0.06375'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.51%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'" + '20.91'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +2.43%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'" + '0.07847'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.21%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'" + '4.530'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +0.15%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'" + '1.644.24'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -1.75%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'" + '1.897.41'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +0.77%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'" + '0.5606'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.97%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.02%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'" + '65.87'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.71%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'" + '26.331.62'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.58%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.01%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'" + '4.729'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +2.80%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'" + '197.88'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +3.38%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'" + '10.28'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.83%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'" + '6.052'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.41%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'" + '1.010'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +0.09%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'" + '146.29'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +1.84%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'" + '0.1222'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.26%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'" + '7.257'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.77%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'" + '16.17'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.08%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'" + '1.506'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +2.32%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'" + '0.05905'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +2.23%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'" + '1.284'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +0.53%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'" + '3.555'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.11%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'" + '3.329'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +1.79%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'" + '1.607'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.88%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'MXToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(35, 4).Value = "'" + '2.830'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.57%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).Value = "'" + '0.9621'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.95%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'" + '2.432'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.31%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'" + '0.5831'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +1.04%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'" + '0.01616'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.89%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'" + '5.956'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +2.40%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'" + '1.079.28'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +3.12%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'" + '0.8579'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +0.76%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.06%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'" + '102.98'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.92%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'" + '1.806.72'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.59%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'" + '58.66'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +3.17%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +1.11%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.79%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.89%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'" + '8.078'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +2.80%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'" + '0.05151'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.11%  '
